$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.511.70'
$ws.Range('E2').Value = '  +2.45%  '

$ws.Range('D3').Value = '2.197.13'
$ws.Range('E3').Value = '  +1.66%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.39'
$ws.Range('E5').Value = '  +5.91%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.612'
$ws.Range('E6').Value = '  +0.95%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.76'
$ws.Range('E7').Value = '  +2.75%  '

$ws.Range('E8').Value = '  -0.16%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.585'
$ws.Range('E9').Value = '  +1.52%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.88'
$ws.Range('E10').Value = '  +0.33%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0916'
$ws.Range('E11').Value = '  +1.19%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.80'
$ws.Range('E12').Value = '  +1.55%  '

$ws.Range('E13').Value = '  +1.48%  '

$ws.Range('D14').Value = '2.528.67'
$ws.Range('E14').Value = '  +1.79%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.31'
$ws.Range('E15').Value = '  +1.85%  '

$ws.Range('D16').Value = '2.182.12'
$ws.Range('E16').Value = '  +1.13%  '

$ws.Range('E17').Value = '  -0.50%  '

$ws.Range('D18').Value = '42.456.95'
$ws.Range('E18').Value = '  +2.68%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000102'
$ws.Range('E19').Value = '  +0.00%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.03'
$ws.Range('E20').Value = '  +2.06%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.92'
$ws.Range('E21').Value = '  +2.63%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.20'
$ws.Range('E22').Value = '  +9.15%  '

$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.65'
$ws.Range('E23').Value = '  -1.65%  '

$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '228.18'
$ws.Range('E24').Value = '  +1.20%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.62'
$ws.Range('E26').Value = '  -0.26%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.36'
$ws.Range('E27').Value = '  +2.30%  '

$ws.Range('E28').Value = '  +1.22%  '

$ws.Range('E29').Value = '  +2.87%  '

$ws.Range('E30').Value = '  +11.33%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '168.87'
$ws.Range('E31').Value = '  -0.51%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.04'
$ws.Range('E32').Value = '  +1.52%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0801'
$ws.Range('E33').Value = '  +4.42%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.16'
$ws.Range('E34').Value = '  +0.39%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.120'
$ws.Range('E35').Value = '  +0.15%  '

$ws.Range('E36').Value = '  +2.45%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.36'
$ws.Range('E37').Value = '  +1.83%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0330'
$ws.Range('E38').Value = '  +10.12%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.06'
$ws.Range('E39').Value = '  -0.29%  '

$ws.Range('E40').Value = '  -0.03%  '

$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.197'
$ws.Range('E41').Value = '  +4.67%  '

$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.24'
$ws.Range('E42').Value = '  -1.30%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '58.99'
$ws.Range('E43').Value = '  +0.18%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.481'
$ws.Range('E44').Value = '  +21.19%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.51'
$ws.Range('E45').Value = '  +6.86%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.33'
$ws.Range('E46').Value = '  -0.18%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0977'
$ws.Range('E47').Value = '  +2.04%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.41'
$ws.Range('E48').Value = '  +11.27%  '

$ws.Range('E49').Value = '  +2.35%  '

$ws.Range('E50').Value = '  +1.95%  '

$ws.Range('E51').Value = '  +1.47%  '
